# Rename the "_merge"/"_diff_days"/"_abs_diff_days"/"_duplicates" result
# columns produced by the (old) date-proximity operator to their new
# "_mp_*" (macpie) prefixed names. This mirrors the project-wide rename
# from "lavalinker" to "macpie" described in the commit message.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("O1").Value = "_mp_merge"
$ws.Range("P1").Value = "_mp_diff_days"
$ws.Range("Q1").Value = "_mp_abs_diff_days"
$ws.Range("R1").Value = "_mp_duplicates"

# The new header labels are longer than the old ones, so the columns'
# "best fit" width grows accordingly (as Excel does when it recomputes
# the best-fit width for a column after the displayed text changes).
$ws.Columns.Item(15).ColumnWidth = 9.5
$ws.Columns.Item(16).ColumnWidth = 11.166666666666666
$ws.Columns.Item(17).ColumnWidth = 14.5
$ws.Columns.Item(18).ColumnWidth = 12.5
